$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking text values remain stored as Text (matches source workbook,
# where every Price/Volume cell is an inline string, even when it looks like a number).

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '72.316.62'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -0.15%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.656.12'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +1.22%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '597.17'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -1.10%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '174.67'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -2.32%  '
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -0.28%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.657.01'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +1.23%  '
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -2.54%  '
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +2.15%  '
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +1.40%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.99'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -0.89%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.144.64'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +1.28%  '
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -1.77%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '72.230.17'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -0.15%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '26.31'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -1.31%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.657.11'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +1.21%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.24'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +5.55%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '8.15'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +1.69%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '370.10'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -2.75%  '
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -0.22%  '
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +1.06%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '72.11'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -1.49%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '4.32'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -1.37%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.75'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -2.21%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.796.28'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +1.26%  '
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -0.04%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0₃0973'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +1.85%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.11'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +0.46%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '499.76'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -3.91%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.29'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -2.60%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.83'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -0.39%  '
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -0.03%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '162.83'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -1.58%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '19.51'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +0.94%  '
$ws.Range('B38').Value = 'Kaspa'
$ws.Range('C38').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.111'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +1.02%  '
$ws.Range('B39').Value = 'WhiteBITCoin'
$ws.Range('C39').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '18.92'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -0.83%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -1.77%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.77'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -3.23%  '
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -0.04%  '
$ws.Range('B43').Value = 'dogwifhat'
$ws.Range('C43').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.59'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +0.29%  '
$ws.Range('B44').Value = 'RenderToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.00'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -1.45%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.332'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +0.26%  '
$ws.Range('B46').Value = 'OKB'
$ws.Range('C46').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '39.50'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -0.04%  '
$ws.Range('B47').Value = 'Aave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '156.13'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +4.25%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.74'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  +0.98%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.558'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +2.63%  '
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +1.79%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0757'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -1.18%  '
